$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "season record" columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the season record for every data row (2 through 44)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 97   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 65   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
